# Apply the "parallel" contingencies update:
#  - extend the table from columns A:O to A:Q (two new columns, P and Q)
#  - row 1 gets new header values 14 (P1) and 15 (Q1), carrying the same
#    style as the rest of the header row
#  - rows 2-25: swap I<->K and M<->O pairs (1<->2) and append P/Q = 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: add P1 / Q1, matching the existing header formatting ---
$ws.Range("N1:O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2   # P (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q (new)
}
